$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.607.46'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '1.584.01'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '213.48'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '24.17'
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0894'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').Value = '1.810.50'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').Value = '1.587.64'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.72'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').Value = '28.616.08'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '62.27'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '231.98'
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.41'
$ws.Range('E20').Value = '  -0.74%  '
$ws.Range('D21').Value = '0.0₃0693'
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  -3.56%  '
$ws.Range('E24').Value = '  -1.48%  '
$ws.Range('E25').Value = '  +5.73%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '151.66'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '15.08'
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.46'
$ws.Range('E28').Value = '  -1.56%  '
$ws.Range('E29').Value = '  -2.03%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +2.43%  '
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('E33').Value = '  -1.04%  '
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('D35').Value = '1.400.16'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.07'
$ws.Range('E36').Value = '  +4.75%  '
$ws.Range('E37').Value = '  -3.86%  '
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('E39').Value = '  +3.21%  '
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.524'
$ws.Range('E41').Value = '  -3.13%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.796'
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('E44').Value = '  +1.64%  '
$ws.Range('E46').Value = '  -2.81%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.961'
$ws.Range('E47').Value = '  -2.15%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '63.35'
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('D49').Value = '1.721.91'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '86.80'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('E51').Value = '  -1.78%  '
